$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.091.73'
$ws.Cells.Item(2, 5).Value = '  +0.38%  '

$ws.Cells.Item(3, 4).Value = '1.835.33'
$ws.Cells.Item(3, 5).Value = '  +0.39%  '

$ws.Cells.Item(4, 5).Value = '  +0.31%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '243.13'
$ws.Cells.Item(5, 5).Value = '  -0.43%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.6154'
$ws.Cells.Item(6, 5).Value = '  -2.41%  '

$ws.Cells.Item(7, 5).Value = '  +0.41%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.07462'
$ws.Cells.Item(8, 5).Value = '  -0.71%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.2921'
$ws.Cells.Item(9, 5).Value = '  -0.50%  '

$ws.Cells.Item(10, 5).Value = '  +1.10%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07687'
$ws.Cells.Item(11, 5).Value = '  -0.24%  '

$ws.Cells.Item(12, 4).Value = '1.833.02'
$ws.Cells.Item(12, 5).Value = '  +0.10%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '5.005'
$ws.Cells.Item(13, 5).Value = '  +0.44%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.6719'
$ws.Cells.Item(14, 5).Value = '  +0.35%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '82.59'
$ws.Cells.Item(15, 5).Value = '  -0.39%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.000009260'
$ws.Cells.Item(16, 5).Value = '  -3.86%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '5.928'
$ws.Cells.Item(17, 5).Value = '  -2.34%  '

$ws.Cells.Item(18, 4).Value = '29.085.86'
$ws.Cells.Item(18, 5).Value = '  +0.28%  '

$ws.Cells.Item(19, 4).Value = '2.086.32'
$ws.Cells.Item(19, 5).Value = '  +0.25%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '231.03'
$ws.Cells.Item(20, 5).Value = '  +1.93%  '

$ws.Cells.Item(21, 5).Value = '  +0.94%  '

$ws.Cells.Item(22, 5).Value = '  +0.50%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '7.180'
$ws.Cells.Item(23, 5).Value = '  +0.39%  '

$ws.Cells.Item(24, 5).Value = '  +0.38%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '160.09'
$ws.Cells.Item(25, 5).Value = '  +0.23%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.1387'

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '8.497'
$ws.Cells.Item(27, 5).Value = '  -0.27%  '

$ws.Cells.Item(28, 5).Value = '  -0.50%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.496'
$ws.Cells.Item(29, 5).Value = '  +0.02%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.158'
$ws.Cells.Item(30, 5).Value = '  +1.05%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '4.131'
$ws.Cells.Item(31, 5).Value = '  +1.72%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.05515'
$ws.Cells.Item(32, 5).Value = '  +2.73%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.211'
$ws.Cells.Item(33, 5).Value = '  +1.36%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.7430'
$ws.Cells.Item(34, 5).Value = '  +0.08%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.838'
$ws.Cells.Item(35, 5).Value = '  -0.95%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.140'
$ws.Cells.Item(36, 5).Value = '  +0.29%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.660'
$ws.Cells.Item(37, 5).Value = '  +0.42%  '

$ws.Cells.Item(38, 5).Value = '  +0.73%  '

$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.01780'
$ws.Cells.Item(39, 5).Value = '  -0.22%  '

$ws.Cells.Item(40, 2).Value = 'Maker'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(40, 4).Value = '1.214.30'
$ws.Cells.Item(40, 5).Value = '  -2.37%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '6.475'
$ws.Cells.Item(41, 5).Value = '  -2.24%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.8946'
$ws.Cells.Item(42, 5).Value = '  -0.70%  '

$ws.Cells.Item(43, 5).Value = '  +0.24%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '102.07'
$ws.Cells.Item(44, 5).Value = '  +0.57%  '

$ws.Cells.Item(45, 4).Value = '1.986.54'
$ws.Cells.Item(45, 5).Value = '  +0.22%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '65.57'
$ws.Cells.Item(46, 5).Value = '  +1.24%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.00000000124'
$ws.Cells.Item(47, 5).Value = '  -0.77%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.5092'
$ws.Cells.Item(48, 5).Value = '  -0.07%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.4066'
$ws.Cells.Item(49, 5).Value = '  -0.04%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '9.101'
$ws.Cells.Item(50, 5).Value = '  +1.56%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.05830'
$ws.Cells.Item(51, 5).Value = '  +1.16%  '
